# File Handling Complete: to CSV | to XLSX
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20 - File Handling: mark finishes, on time, with real finish date
$ws.Range("D20").Value = 44456
$ws.Range("E20").Value = "Finishes"
$ws.Range("F20").Value = "OnTime"

# Row 13 - Admissions Process (to CSV): mark finished, delayed, with real finish date
$ws.Range("D13").Value = 44456
$ws.Range("E13").Value = "Finished"
$ws.Range("F13").Value = "Delayed "

# Row 14 - Search Operations (to XLSX): mark finished, delayed, with real finish date
$ws.Range("D14").Value = 44456
$ws.Range("E14").Value = "Finished"
$ws.Range("F14").Value = "Delayed "

# Row 15 - Count Operations: mark finished, delayed, with real finish date
$ws.Range("D15").Value = 44456
$ws.Range("E15").Value = "Finished"
$ws.Range("F15").Value = "Delayed "

# Row 21 - DataBase Operations: mark early
$ws.Range("F21").Value = "Early "

# Scroll the view down and move the selection to reflect where work left off
$ws.Application.Goto($ws.Range("A12"), $true)
$ws.Range("F22").Select()
